$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: replace the "A" ticket entry with the "John Lennon" ticket entry
$ws.Range("C3").Value = "John Lennon"
$ws.Range("D3").Value = "05.01.2022"
$ws.Range("E3").Value = """Analysis for Excel"" Add-In causes crashes"
$ws.Range("F3").Value = "Hello, Is there any more stable version of ""Analysis for Excel"" for Office 365? This one causes a lot of crashes when run with VBA. My current version is as attached - 2.8."
$ws.Range("G3").Value = '[{"start":45, "end":60, "key":"Fehlerbeschreibung"},{"start":30, "end":45, "key":"System"}]'
